$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at U:V. This shifts the existing "04-30_0" column
# (U) two places to the right, landing it at W, while leaving everything to
# the left (through T) untouched. The two freshly-inserted columns (U, V)
# start out empty, inheriting the neighboring column's formatting.
$ws.Columns("U:V").Insert()

# New header cells for the extra "05-01" day: _A (attempts/status) then _0
# (running total). The inserted columns already inherited the header style
# (s="1") from their neighbors, matching every other header cell.
$ws.Range("U1").Value = "05-01_A"
$ws.Range("V1").Value = "05-01_0"

for ($r = 2; $r -le 119; $r++) {
    # New "05-01_0" value equals the old "04-30_0" value (now shifted to W),
    # but re-typed as a real number instead of the legacy inline-string.
    $wVal = $ws.Cells.Item($r, 23).Value()
    $ws.Cells.Item($r, 21).Value = $wVal
    $ws.Cells.Item($r, 21).Style = "Normal"

    # New "05-01_A" column carries forward the same status/style as "04-30_A".
    $tVal = $ws.Cells.Item($r, 20).Value()
    $ws.Cells.Item($r, 22).Value = $tVal
}
